$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.239.84'
$ws.Range("E2").Value = '  -1.26%  '

$ws.Range("D3").Value = '1.866.92'
$ws.Range("E3").Value = '  -1.50%  '

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = '  +0.36%  '

$ws.Range("D5").Value = "'236.92"
$ws.Range("E5").Value = '  -1.14%  '

$ws.Range("E6").Value = '  +0.49%  '

$ws.Range("D7").Value = "'0.4784"
$ws.Range("E7").Value = '  -2.90%  '

$ws.Range("D8").Value = "'0.2830"
$ws.Range("E8").Value = '  -4.02%  '

$ws.Range("D9").Value = "'0.06521"
$ws.Range("E9").Value = '  -3.08%  '

$ws.Range("D10").Value = '1.882.37'
$ws.Range("E10").Value = '  -0.62%  '

$ws.Range("D11").Value = "'0.07377"
$ws.Range("E11").Value = '  +0.25%  '

$ws.Range("D12").Value = "'16.46"
$ws.Range("E12").Value = '  -4.20%  '

$ws.Range("D13").Value = "'5.163"
$ws.Range("E13").Value = '  +0.21%  '

$ws.Range("D14").Value = "'87.29"
$ws.Range("E14").Value = '  -1.27%  '

$ws.Range("D15").Value = "'0.6489"
$ws.Range("E15").Value = '  -3.27%  '

$ws.Range("D16").Value = '30.198.87'

$ws.Range("D17").Value = "'13.29"
$ws.Range("E17").Value = '  -1.33%  '

$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = '  +0.10%  '

$ws.Range("D19").Value = "'0.000007611"
$ws.Range("E19").Value = '  -3.48%  '

$ws.Range("D20").Value = '2.133.81'
$ws.Range("E20").Value = '  +0.07%  '

$ws.Range("D21").Value = "'1.004"
$ws.Range("E21").Value = '  +0.29%  '

$ws.Range("D22").Value = "'5.278"
$ws.Range("E22").Value = '  -0.52%  '

$ws.Range("D23").Value = "'216.84"
$ws.Range("E23").Value = '  +13.09%  '

$ws.Range("D24").Value = "'6.127"
$ws.Range("E24").Value = '  -1.56%  '

$ws.Range("D25").Value = "'9.318"
$ws.Range("E25").Value = '  -2.55%  '

$ws.Range("D26").Value = "'164.34"
$ws.Range("E26").Value = '  +2.15%  '

$ws.Range("D27").Value = "'18.58"
$ws.Range("E27").Value = '  +0.63%  '

$ws.Range("D28").Value = "'1.908"
$ws.Range("E28").Value = '  -1.95%  '

$ws.Range("E29").Value = '  -2.06%  '

$ws.Range("D30").Value = "'4.263"
$ws.Range("E30").Value = '  -4.26%  '

$ws.Range("D31").Value = "'0.09170"
$ws.Range("E31").Value = '  -0.29%  '

$ws.Range("D32").Value = "'3.979"
$ws.Range("E32").Value = '  -5.15%  '

$ws.Range("D33").Value = "'0.05035"
$ws.Range("E33").Value = '  -3.97%  '

$ws.Range("D34").Value = "'0.7451"
$ws.Range("E34").Value = '  -0.33%  '

$ws.Range("D35").Value = "'1.134"
$ws.Range("E35").Value = '  +2.41%  '

$ws.Range("D36").Value = "'2.692"
$ws.Range("E36").Value = '  -0.78%  '

$ws.Range("D37").Value = "'0.01834"
$ws.Range("E37").Value = '  -0.18%  '

$ws.Range("D38").Value = "'2.623"
$ws.Range("E38").Value = '  -2.67%  '

$ws.Range("D39").Value = "'0.9096"
$ws.Range("E39").Value = '  -1.50%  '

$ws.Range("D40").Value = "'2.064"
$ws.Range("E40").Value = '  +0.16%  '

$ws.Range("D41").Value = "'5.930"
$ws.Range("E41").Value = '  -0.53%  '

$ws.Range("D42").Value = "'106.67"
$ws.Range("E42").Value = '  +0.35%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = "'0.4268"
$ws.Range("E43").Value = '  -3.64%  '

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = "'1.006"
$ws.Range("E44").Value = '  +1.15%  '

$ws.Range("D45").Value = "'7.463"
$ws.Range("E45").Value = '  -1.79%  '

$ws.Range("D46").Value = "'1.569"
$ws.Range("E46").Value = '  +10.08%  '

$ws.Range("D47").Value = "'0.1309"
$ws.Range("E47").Value = '  -5.74%  '

$ws.Range("D48").Value = "'64.37"
$ws.Range("E48").Value = '  -9.66%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = "'34.39"
$ws.Range("E49").Value = '  -2.02%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = "'8.781"
$ws.Range("E50").Value = '  -3.30%  '

$ws.Range("D51").Value = "'0.05714"
$ws.Range("E51").Value = '  -2.01%  '

